# Weekly update: prepend the newest week's data (2 new rows) to the top of
# the data block (right after the header row + the existing row that stays
# in place), pushing all older rows down by two rows. No rows are removed;
# the used range grows from A1:R118 to A1:R120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 11:12 - this shifts old rows 11..118 down to 13..120
# and (as Excel does by default) carries the formatting of the row above
# into the freshly inserted rows.
$ws.Rows("11:12").Insert()

# New data for the newest reporting week.
$row11 = @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44532, 4, 100112031, "Poroto verde", "Magnum", "Primera", 700, 14000, 16000, 15000, "$/malla 25 kilos", "Provincia de Limarí", 600, 25, "Hortaliza")
$row12 = @(2, "Comercializadora del Agro de Limarí", "Coquimbo", 44532, 4, 100112031, "Poroto verde", "Sin especificar", "Primera", 500, 25000, 27000, 26000, "$/malla 25 kilos", "Provincia de Limarí", 1040, 25, "Hortaliza")

for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, $i + 1).Value = $row11[$i]
}
for ($i = 0; $i -lt $row12.Length; $i++) {
    $ws.Cells.Item(12, $i + 1).Value = $row12[$i]
}
